$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-CellText 'D2' '67.470.75'
Set-CellText 'E2' '  -0.95%  '

Set-CellText 'D3' '3.342.13'
Set-CellText 'E3' '  +2.34%  '

Set-CellText 'E4' '  +0.02%  '

Set-CellText 'D5' '579.39'
Set-CellText 'E5' '  -0.63%  '

Set-CellText 'D6' '183.34'
Set-CellText 'E6' '  -0.74%  '

Set-CellText 'E7' '  +0.05%  '

Set-CellText 'D8' '0.604'
Set-CellText 'E8' '  +0.65%  '

Set-CellText 'D9' '0.129'
Set-CellText 'E9' '  -0.39%  '

Set-CellText 'D10' '6.67'
Set-CellText 'E10' '  +0.85%  '

Set-CellText 'D11' '0.407'
Set-CellText 'E11' '  +0.05%  '

Set-CellText 'D12' '3.923.02'
Set-CellText 'E12' '  +2.28%  '

Set-CellText 'E13' '  -0.64%  '

Set-CellText 'D14' '27.37'
Set-CellText 'E14' '  +0.28%  '

Set-CellText 'D15' '67.668.35'
Set-CellText 'E15' '  -0.62%  '

Set-CellText 'E16' '  -0.14%  '

Set-CellText 'D17' '3.330.70'
Set-CellText 'E17' '  +1.76%  '

Set-CellText 'E18' '  +6.98%  '

Set-CellText 'E19' '  +2.37%  '

Set-CellText 'D20' '5.67'
Set-CellText 'E20' '  -0.70%  '

Set-CellText 'D21' '7.73'
Set-CellText 'E21' '  +2.79%  '

Set-CellText 'D22' '73.97'
Set-CellText 'E22' '  +3.75%  '

Set-CellText 'E23' '  -0.15%  '

Set-CellText 'D24' '3.476.89'
Set-CellText 'E24' '  +2.02%  '

Set-CellText 'D25' '0.512'
Set-CellText 'E25' '  +1.03%  '

Set-CellText 'E26' '  +1.72%  '

Set-CellText 'D27' '0.192'
Set-CellText 'E27' '  +3.11%  '

Set-CellText 'D28' '9.07'
Set-CellText 'E28' '  -3.01%  '

Set-CellText 'D29' '1.00'
Set-CellText 'E29' '  -0.06%  '

Set-CellText 'E30' '  +1.02%  '

Set-CellText 'D31' '22.94'
Set-CellText 'E31' '  +1.65%  '

Set-CellText 'D32' '5.34'
Set-CellText 'E32' '  -2.00%  '

Set-CellText 'D34' '6.79'
Set-CellText 'E34' '  -0.48%  '

Set-CellText 'E35' '  -0.55%  '

Set-CellText 'D36' '1.50'
Set-CellText 'E36' '  +4.37%  '

Set-CellText 'D37' '161.91'
Set-CellText 'E37' '  -1.23%  '

Set-CellText 'B38' 'EnergySwap'
Set-CellText 'C38' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText 'D38' '27.36'
Set-CellText 'E38' '  +2.21%  '

Set-CellText 'B39' 'Stacks'
Set-CellText 'C39' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-CellText 'D39' '1.85'
Set-CellText 'E39' '  -1.62%  '

Set-CellText 'D40' '2.823.58'
Set-CellText 'E40' '  +7.23%  '

Set-CellText 'D41' '0.791'
Set-CellText 'E41' '  -0.53%  '

Set-CellText 'D42' '4.46'
Set-CellText 'E42' '  +0.76%  '

Set-CellText 'D43' '6.22'
Set-CellText 'E43' '  -1.47%  '

Set-CellText 'D44' '40.36'
Set-CellText 'E44' '  -0.65%  '

Set-CellText 'D45' '0.0673'
Set-CellText 'E45' '  +0.11%  '

Set-CellText 'D46' '24.67'
Set-CellText 'E46' '  +2.19%  '

Set-CellText 'E47' '  -2.30%  '

Set-CellText 'D48' '324.79'
Set-CellText 'E48' '  -3.57%  '

Set-CellText 'D49' '0.0273'
Set-CellText 'E49' '  +0.14%  '

Set-CellText 'D50' '0.987'
Set-CellText 'E50' '  +0.80%  '

Set-CellText 'D51' '31.05'
Set-CellText 'E51' '  +1.77%  '
